$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 66
$ws.Range("F2").Value = 46
$ws.Range("H2").Value = 58
$ws.Range("F3").Value = 18
$ws.Range("H3").Value = 22
$ws.Range("F5").Value = 6
$ws.Range("H5").Value = 6
$ws.Range("E6").Value = 9
$ws.Range("F9").Value = 17
$ws.Range("H9").Value = 26
$ws.Range("F12").Value = 5
$ws.Range("H12").Value = 5
$ws.Range("E15").Value = 183
$ws.Range("F15").Value = 107
$ws.Range("H15").Value = 148
$ws.Range("E17").Value = 143
$ws.Range("F17").Value = 76
$ws.Range("H17").Value = 108
$ws.Range("E18").Value = 138
$ws.Range("F18").Value = 68
$ws.Range("H18").Value = 105
$ws.Range("E19").Value = 72
$ws.Range("F19").Value = 45
$ws.Range("H19").Value = 58
$ws.Range("F23").Value = 5
$ws.Range("H23").Value = 7
$ws.Range("F24").Value = 19
$ws.Range("H24").Value = 23
$ws.Range("E25").Value = 27
$ws.Range("F25").Value = 17
$ws.Range("H25").Value = 25
$ws.Range("F26").Value = 22
$ws.Range("H26").Value = 32
$ws.Range("F27").Value = 13
$ws.Range("H27").Value = 17
$ws.Range("F32").Value = 8
$ws.Range("H32").Value = 17
$ws.Range("F33").Value = 16
$ws.Range("H33").Value = 28
$ws.Range("E34").Value = 26
$ws.Range("F34").Value = 11
$ws.Range("H34").Value = 14
$ws.Range("F35").Value = 8
$ws.Range("H35").Value = 9
$ws.Range("F36").Value = 65
$ws.Range("H36").Value = 97
$ws.Range("F37").Value = 40
$ws.Range("H37").Value = 52
$ws.Range("E38").Value = 90
$ws.Range("F38").Value = 24
$ws.Range("H38").Value = 44
$ws.Range("F39").Value = 18
$ws.Range("H39").Value = 26
$ws.Range("F40").Value = 20
$ws.Range("H40").Value = 22
$ws.Range("F41").Value = 27
$ws.Range("H41").Value = 38
$ws.Range("F42").Value = 28
$ws.Range("H42").Value = 37
$ws.Range("F43").Value = 22
$ws.Range("H43").Value = 25
$ws.Range("E44").Value = 35
$ws.Range("F44").Value = 21
$ws.Range("H44").Value = 31
$ws.Range("E45").Value = 29
$ws.Range("E47").Value = 69
$ws.Range("F47").Value = 44
$ws.Range("H47").Value = 54
$ws.Range("E48").Value = 46
$ws.Range("F48").Value = 31
$ws.Range("H48").Value = 37
$ws.Range("E49").Value = 84
$ws.Range("F49").Value = 47
$ws.Range("H49").Value = 64
$ws.Range("F50").Value = 13
$ws.Range("H50").Value = 22
$ws.Range("F53").Value = 4
$ws.Range("H53").Value = 6
$ws.Range("F58").Value = 5
$ws.Range("H58").Value = 5
$ws.Range("F59").Value = 6
$ws.Range("H59").Value = 10
$ws.Range("F60").Value = 14
$ws.Range("H60").Value = 19
$ws.Range("F61").Value = 18
$ws.Range("H61").Value = 28
$ws.Range("F62").Value = 17
$ws.Range("H62").Value = 31
$ws.Range("E64").Value = 41
$ws.Range("F65").Value = 16
$ws.Range("H65").Value = 29
$ws.Range("F66").Value = 29
$ws.Range("H66").Value = 37
$ws.Range("F69").Value = 10
$ws.Range("H69").Value = 13
$ws.Range("E70").Value = 53
$ws.Range("F70").Value = 29
$ws.Range("H70").Value = 42
$ws.Range("F71").Value = 24
$ws.Range("H71").Value = 34
$ws.Range("E73").Value = 36
$ws.Range("F73").Value = 17
$ws.Range("H73").Value = 29
$ws.Range("F75").Value = 11
$ws.Range("H75").Value = 16
$ws.Range("F76").Value = 24
$ws.Range("H76").Value = 41
$ws.Range("F77").Value = 26
$ws.Range("H77").Value = 43
$ws.Range("F79").Value = 26
$ws.Range("H79").Value = 37
$ws.Range("F80").Value = 19
$ws.Range("H80").Value = 31
$ws.Range("E85").Value = 7
$ws.Range("E87").Value = 22
$ws.Range("F87").Value = 9
$ws.Range("H87").Value = 16
$ws.Range("F88").Value = 23
$ws.Range("H88").Value = 31
$ws.Range("E89").Value = 50
